# Add files via upload
#
# On the "axes" sheet, insert three new columns (D:F) carrying an
# "*_arrow" header row and a "<component> (w/w)" data row, pushing the
# existing Title/"USDA Textural Sediment Classification" column from D to G.
# Also (re)select the "axes" sheet/cell as the active view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("axes")

# Shift the existing "Title" column (currently D) three columns to the right
# (to G) by inserting three new whole columns at D:F.
$ws.Range("D1:F2").EntireColumn.Insert()

# New header row (row 1)
$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

# New data row (row 2)
$ws.Range("D2").Value = "Clay (w/w)"
$ws.Range("E2").Value = "Sand (w/w)"
$ws.Range("F2").Value = "Silt (w/w)"

# Widen the new columns
$ws.Range("D1:F2").ColumnWidth = 13.3

# Make "axes" the active sheet/tab and select E5 (matches the saved view)
$ws.Activate()
$ws.Range("E5").Select()
